# Update the date line and the multiplication problems in the table.
# Uses Find/Execute with MatchWholeWord=$true and MatchWildcards=$false
# to replace each exact text string once.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-02-29 Thursday"; New = "2024-03-01 Friday" },
    @{ Old = "43×88="; New = "98×84=" },
    @{ Old = "60×73="; New = "21×97=" },
    @{ Old = "96×41="; New = "91×97=" },
    @{ Old = "18×49="; New = "61×27=" },
    @{ Old = "85×70="; New = "27×19=" },
    @{ Old = "66×77="; New = "60×87=" },
    @{ Old = "86×62="; New = "20×95=" },
    @{ Old = "72×81="; New = "93×13=" },
    @{ Old = "27×22="; New = "67×22=" },
    @{ Old = "74×84="; New = "19×70=" },
    @{ Old = "65×66="; New = "60×37=" },
    @{ Old = "18×45="; New = "26×27=" },
    @{ Old = "42×45="; New = "28×39=" },
    @{ Old = "54×48="; New = "64×29=" },
    @{ Old = "63×92="; New = "55×79=" },
    @{ Old = "86×73="; New = "18×19=" },
    @{ Old = "47×37="; New = "17×68=" },
    @{ Old = "46×35="; New = "14×22=" },
    @{ Old = "85×13="; New = "94×96=" },
    @{ Old = "66×53="; New = "43×40=" },
    @{ Old = "25×15="; New = "12×31=" },
    @{ Old = "74×65="; New = "67×47=" },
    @{ Old = "36×81="; New = "67×52=" },
    @{ Old = "28×75="; New = "41×77=" },
    @{ Old = "79×50="; New = "12×51=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
